# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row (row 46 in this workbook) and the header row (row 1).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Header cells - copy style/formatting from the existing header (column A) so
# the new headers look like the rest of the header row.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows: every team row gets the same season record.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 86   # AD = Wins
    $ws.Cells.Item($r, 31).Value = 76   # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = Ties
}
